# Apply the "twiddling numbers - union plaza" edit:
#  - new a/b linear-fit coefficients for longitude (rows 76-77) and
#    latitude (rows 81-82), with labels in column G
#  - new "Latitudes"/"Longitudes" labels (rows 74 & 79)
#  - a new little X[mm]/Y[mm]/X[in]/Y[in] header row (85) followed by
#    36 rows (86-121) of projected point coordinates in mm and inches
#  - scroll/selection moved down to the newly active working area

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- labels -----------------------------------------------------------
$ws.Range("G74").Value = "Longitudes"

$ws.Range("G76").Value = "a="
$ws.Range("G77").Value = "b="

$ws.Range("G79").Value = "Latitudes"

$ws.Range("G81").Value = "a="
$ws.Range("G82").Value = "b="

# ---- linear fit coefficients -------------------------------------------
# Longitude (X) fit: maps GPS longitude (G61/G66) -> mm position (H61/H66)
$ws.Range("H76").Formula = "=(H61-H66)/(G61-G66)"
$ws.Range("H77").Formula = "=H61-G61*H76"

# Latitude (Y) fit: maps GPS latitude (G60/G65) -> mm position (H60/H65)
$ws.Range("H81").Formula = "=(H60-H65)/(G60-G65)"
$ws.Range("H82").Formula = "=H61-G61*H81"

# ---- small header row for the projected-point table --------------------
$ws.Range("H85").Value = "X [mm]"
$ws.Range("I85").Value = "Y [mm]"
$ws.Range("K85").Value = "X [in]"
$ws.Range("L85").Value = "Y [in]"

# ---- 36 rows (86-121) of projected coordinates --------------------------
# Row r uses source lat/long pair at row (r-54): 86->32 ... 121->67
for ($r = 86; $r -le 121; $r++) {
    $src = $r - 54
    $ws.Cells.Item($r, 8).Formula  = "=`$H`$76*C$src+`$H`$77"   # H: X [mm]
    $ws.Cells.Item($r, 9).Formula  = "=B$src*`$H`$81+`$H`$82"   # I: Y [mm]
    $ws.Cells.Item($r, 11).Formula = "=H$r/25.4"                # K: X [in]
    $ws.Cells.Item($r, 12).Formula = "=I$r/25.4"                # L: Y [in]
}

# ---- recalc so cached <v> values are written out ------------------------
$ws.Calculate()

# ---- view state: active selection (engine doesn't persist topLeftCell) --
$ws.Range("F61").Select()
